$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A5").Value = "Beitragsbemessungsgrenze GKV"
$ws.Range("A6").Value = "Jahresarbeitsentgeltgrenze GKV"

$ws.Range("A4").Select()
